$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full country list sorted by "Casos totales" (column B) descending, reflecting the
# updated data pull. Each entry: Country, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @("China",81218,47,73650,4287,1399,4,3281),
    @("Italia",74386,5210,9362,57521,3489,683,7503),
    @("Estados Unidos",60653,5797,387,59447,1307,39,819),
    @("España",47610,5552,5367,38809,2636,443,3434),
    @("Alemania",35740,2749,3540,32014,23,27,186),
    @("Iran",27017,2206,9625,15315,0,143,2077),
    @("Francia",22304,0,3281,17923,2516,0,1100),
    @("Suiza",10537,660,131,10257,141,27,149),
    @("Corea del Sur",9137,100,3730,5281,59,6,126),
    @("Reino Unido",8227,150,135,7659,20,11,433),
    @("Paises Bajos",6412,852,3,6053,582,80,356),
    @("Austria",5588,305,9,5549,28,2,30),
    @("Belgica",4937,668,547,4212,474,56,178),
    @("Noruega",3034,168,6,3014,57,2,14),
    @("Portugal",2995,633,22,2930,61,10,43),
    @("Canada",2923,131,185,2710,1,2,28),
    @("Suecia",2526,227,16,2466,158,4,44),
    @("Australia",2431,114,118,2304,11,1,9),
    @("Brasil",2274,27,2,2225,18,1,47),
    @("Israel",2170,240,58,2107,37,2,5),
    @("Turquia",1872,0,26,1802,0,0,44),
    @("Malasia",1796,172,199,1577,64,4,20),
    @("Dinamarca",1715,124,1,1680,69,2,34),
    @("Chequia",1654,260,10,1638,19,3,6),
    @("Luxemburgo",1333,234,6,1319,3,0,8),
    @("Irlanda",1329,0,5,1317,29,0,7),
    @("Japon",1193,0,285,865,54,0,43),
    @("Ecuador",1173,91,3,1142,2,1,28),
    @("Chile",1142,220,22,1117,7,1,3),
    @("Polonia",1031,130,2,1015,3,4,14),
    @("Pakistan",1022,50,21,993,5,1,8),
    @("Tailandia",934,107,70,860,11,0,4),
    @("Rumania",906,112,86,806,18,2,14),
    @("Arabia Saudita",900,133,29,869,0,1,2),
    @("Finlandia",880,88,10,867,22,2,3),
    @("Grecia",821,78,36,763,53,2,22),
    @("Indonesia",790,104,31,701,0,3,58),
    @("Islandia",737,89,56,679,11,0,2),
    @("Crucero",712,0,587,115,15,0,10),
    @("Sudafrica",709,155,12,697,2,0,0),
    @("Rusia",658,163,29,628,8,0,1),
    @("Filipinas",636,84,26,572,1,3,38),
    @("Singapur",631,73,160,469,17,0,2),
    @("India",606,70,42,554,0,0,10),
    @("Eslovenia",528,48,10,513,14,1,5),
    @("Catar",526,0,41,485,6,0,0),
    @("Panama",443,0,1,434,33,2,8),
    @("Croacia",442,60,22,419,6,0,1),
    @("Egipto",442,40,93,328,0,1,21),
    @("Barein",419,27,177,239,2,0,3),
    @("Peru",416,0,1,406,9,2,9),
    @("Hong Kong",410,23,102,304,4,0,4),
    @("Mexico",405,38,4,396,1,1,5),
    @("Estonia",404,35,8,395,5,1,1),
    @("Republica Dominicana",392,80,3,379,0,4,10),
    @("Argentina",387,0,52,327,0,2,8),
    @("Serbia",384,81,15,365,21,1,4),
    @("Colombia",378,0,6,369,0,0,3),
    @("Irak",346,30,103,214,0,2,29),
    @("Libano",333,15,8,321,4,0,4),
    @("Emiratos Arabes Unidos",333,85,52,279,2,0,2),
    @("Argelia",302,38,65,216,0,2,21),
    @("Armenia",265,16,16,249,6,0,0),
    @("Lituania",255,46,1,250,1,2,4),
    @("Bulgaria",242,24,4,235,8,0,3),
    @("Taiwan",235,19,29,204,0,0,2),
    @("Hungria",226,39,21,195,6,1,10),
    @("Letonia",221,24,1,220,0,0,0),
    @("Eslovaquia",216,12,7,209,2,0,0),
    @("Nueva Zelanda",205,0,22,183,0,0,0),
    @("Kuwait",195,4,43,152,6,0,0),
    @("Uruguay",189,0,0,189,3,0,0),
    @("Principado de Andorra",188,24,1,186,6,0,1),
    @("San Marino",187,0,4,162,12,0,21),
    @("Costa Rica",177,0,2,173,4,0,2),
    @("Republica de Macedonia",177,29,1,173,1,1,3),
    @("Tunez",173,59,2,166,11,1,5),
    @("Marruecos",170,0,6,159,1,0,5),
    @("Bosnia y Herzegovina",168,0,2,163,1,0,3),
    @("Jordania",154,0,1,153,0,0,0),
    @("Moldavia",149,24,2,146,20,0,1),
    @("Albania",146,23,17,124,3,0,5),
    @("Vietnam",141,7,17,124,3,0,0),
    @("Republica de Chipre",132,8,3,126,3,0,3),
    @("Islas Feroe",132,10,38,94,2,0,0),
    @("Malta",129,19,2,127,1,0,0),
    @("Ucrania",116,14,1,111,0,1,4),
    @("Burkina Faso",114,0,7,103,0,0,4),
    @("Reunion",111,17,1,110,0,0,0),
    @("Brunei",109,5,2,107,1,0,0),
    @("Sri Lanka",102,0,3,99,2,0,0),
    @("Senegal",99,13,9,90,0,0,0),
    @("Oman",99,15,17,82,0,0,0),
    @("Camboya",96,5,10,86,1,0,0),
    @("Azerbaiyan",93,6,10,81,6,1,2),
    @("Venezuela",91,7,15,76,2,0,0),
    @("Bielorrusia",86,5,29,57,2,0,0),
    @("Afganistan",84,10,2,80,0,1,2),
    @("Kazajistan",80,8,0,80,0,0,0),
    @("Guadalupe",73,0,0,72,4,0,1),
    @("Costa de Marfil",73,0,2,71,0,0,0),
    @("Georgia",73,3,10,63,1,0,0),
    @("Camerun",70,4,2,67,0,0,1),
    @("Ghana",68,15,0,65,0,1,3),
    @("Estado de Palestina",64,4,16,47,0,1,1),
    @("Trinidad yTobago",60,3,0,60,0,0,0),
    @("Uzbekistan",59,9,0,59,4,0,0),
    @("Martinica",57,0,0,56,7,0,1),
    @("Cuba",57,9,1,55,2,0,1),
    @("Montenegro",52,5,0,51,0,0,1),
    @("Liechtenstein",51,0,0,51,0,0,0),
    @("Consejo Danes para los Refugiados",48,3,0,46,0,0,2),
    @("Mauricio",48,6,0,46,1,0,2),
    @("Nigeria",46,2,2,43,0,0,1),
    @("Kirguistan",44,2,0,44,0,0,0),
    @("Ruanda",40,0,0,40,0,0,0),
    @("Puerto Rico",39,0,1,36,0,0,2),
    @("Banglades",39,0,7,27,0,1,5),
    @("Paraguay",37,10,0,34,1,1,3),
    @("Mayotte",36,0,0,36,0,0,0),
    @("Honduras",36,6,0,36,0,0,0),
    @("Bolivia",32,3,0,32,0,0,0),
    @("Guam",32,0,0,31,0,0,1),
    @("Macao",30,1,10,20,0,0,0),
    @("Kenia",28,3,1,27,0,0,0),
    @("Monaco",27,4,1,26,0,0,0),
    @("Gibraltar",26,11,5,21,0,0,0),
    @("Polinesia Francesa",25,0,0,25,0,0,0),
    @("Jamaica",25,4,2,22,0,0,1),
    @("Isla de Man",23,0,0,23,0,0,0),
    @("Togo",23,3,1,22,0,0,0),
    @("Guayana Francesa",23,0,6,17,0,0,0),
    @("Guatemala",21,0,0,20,0,0,1),
    @("Madagascar",19,2,0,19,0,0,0),
    @("Barbados",18,0,0,18,0,0,0),
    @("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
    @("Aruba",17,0,1,16,0,0,0),
    @("Uganda",14,5,0,14,0,0,0),
    @("Nueva Caledonia",14,4,0,14,0,0,0),
    @("Maldivas",13,0,8,5,0,0,0),
    @("Etiopia",12,0,0,12,0,0,0),
    @("Zambia",12,9,0,12,0,0,0),
    @("Tanzania",12,0,0,12,0,0,0),
    @("Republica de Yibuti",11,8,0,11,0,0,0),
    @("Mongolia",10,0,0,10,0,0,0),
    @("El Salvador",9,4,0,9,0,0,0),
    @("Guinea Ecuatorial",9,0,0,9,0,0,0),
    @("San Martin (Parte Francesa)",8,0,0,8,0,0,0),
    @("Surinam",8,1,0,8,0,0,0),
    @("Haiti",7,0,0,7,0,0,0),
    @("Seychelles",7,0,0,7,0,0,0),
    @("Dominica",7,0,0,7,0,0,0),
    @("Niger",7,4,0,6,0,1,1),
    @("Namibia",7,0,2,5,0,0,0),
    @("Bermudas",6,0,0,6,0,0,0),
    @("Benin",6,0,0,6,0,0,0),
    @("Gabon",6,0,0,5,0,0,1),
    @("Islas Caimanes",6,0,0,5,0,0,1),
    @("Curazao",6,0,0,5,0,0,1),
    @("Groenlandia",6,1,2,4,0,0,0),
    @("Mozambique",5,2,0,5,0,0,0),
    @("Fiyi",5,1,0,5,0,0,0),
    @("Guyana",5,0,0,4,0,0,1),
    @("Bahamas",5,0,1,4,0,0,0),
    @("Siria",4,3,0,4,0,0,0),
    @("Congo",4,0,0,4,0,0,0),
    @("Santa Sede",4,0,0,4,0,0,0),
    @("Suazilandia",4,0,0,4,0,0,0),
    @("Guinea",4,0,0,4,0,0,0),
    @("Cabo Verde",4,1,0,3,0,0,1),
    @("Liberia",3,0,0,3,0,0,0),
    @("Republica de Africa Central",3,0,0,3,0,0,0),
    @("Santa Lucia",3,0,0,3,0,0,0),
    @("Angola",3,0,0,3,0,0,0),
    @("San Bartolome",3,0,0,3,0,0,0),
    @("San Martin (Parte Holandesa)",3,1,0,3,0,0,0),
    @("Antigua y Barbuda",3,0,0,3,0,0,0),
    @("Republica del Chad",3,0,0,3,0,0,0),
    @("Laos",3,1,0,3,0,0,0),
    @("Birmania",3,0,0,3,0,0,0),
    @("Sudan",3,0,0,2,0,0,1),
    @("Zimbabue",3,0,0,2,0,0,1),
    @("Gambia",3,0,0,2,0,0,1),
    @("Nepal",3,1,1,2,0,0,0),
    @("Butan",2,0,0,2,0,0,0),
    @("Nicaragua",2,0,0,2,0,0,0),
    @("Mali",2,2,0,2,0,0,0),
    @("Guinea-Bisau",2,2,0,2,0,0,0),
    @("Mauritania",2,0,0,2,0,0,0),
    @("Belice",1,0,0,1,0,0,0),
    @("San Vicente y las Granadinas",1,0,0,1,0,0,0),
    @("Libia",1,0,0,1,0,0,0),
    @("Somalia",1,0,0,1,0,0,0),
    @("Timor Oriental",1,0,0,1,0,0,0),
    @("Eritrea",1,0,0,1,0,0,0),
    @("Montserrat",1,0,0,1,0,0,0),
    @("Granada",1,0,0,1,0,0,0),
    @("Papua Nueva Guinea",1,0,0,1,0,0,0),
    @("Islas Turcas y Caicos",1,0,0,1,0,0,0)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
    $ws.Cells.Item($row, 8).Value = $item[7]
    $row = $row + 1
}

# Update the "last updated" footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 18:16"
